# Apply edits to TestData.xlsx (LoginData sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Activate()

# Swap the "Execute" flag: row 3 (L002) becomes YES, row 5 (L004) becomes NO
$ws.Range("A3").Value = "YES"
$ws.Range("A5").Value = "NO"

# Update the active cell selection to D10
$ws.Range("D10").Select()
